$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '91.875.54'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '3.125.25'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '618.37'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.09%  '
$ws.Range("E7").Value = '  -3.52%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.386'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.59%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("D10").Value = '3.119.90'
$ws.Range("E10").Value = '  +0.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.752'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.203'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.69%  '
$ws.Range("E13").Value = '  -0.98%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.63'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.65%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '34.97'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.72%  '
$ws.Range("D16").Value = '91.524.19'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = '3.695.20'
$ws.Range("E17").Value = '  -0.14%  '
$ws.Range("D18").Value = '3.156.45'
$ws.Range("E18").Value = '  +0.68%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.92'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.30%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.82'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '450.51'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.76%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.39'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.84%  '
$ws.Range("B24").Value = 'PEPE'
$ws.Range("C24").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000203'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -8.41%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.88'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '89.75'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.40%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.77'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.141'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +20.31%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.229'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.168'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -9.20%  '
$ws.Range("E33").Value = '  +3.06%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.38'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.81%  '
$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.176'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.82%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.71'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '26.34'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.10%  '
$ws.Range("E38").Value = '  +2.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.93'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.97%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '490.60'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.20%  '
$ws.Range("E41").Value = '  +1.20%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.436'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.67%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.44'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.56%  '
$ws.Range("E44").Value = '  +0.26%  '
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '158.25'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.702'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.14%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.91'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.35'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.11%  '
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '44.19'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.31%  '
$ws.Range("B51").Value = 'Filecoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.41'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.76%  '
